# "Add Magic Details, Began Tutorial"
#
# 1) Spell sheet: add an "Element Type" column (H), shifting the old
#    Effect/Extra/Note columns right (H->J, J->L, L->N), and fully fill in
#    the previously-stub "Undead Screech" row (25).
# 2) Physical Armour sheet: reword the magic-damage armour note (3% -> 20%).
# 3) Selection / active-tab bookkeeping to match where the author ended up.

$wb = $excel.ActiveWorkbook
$spell = $wb.Worksheets.Item("Spell")
$phys  = $wb.Worksheets.Item("Physical Armour")

# --- Finish filling out the new "Undead Screech" row first -----------------
# (matches the authoring order the new shared strings were appended in)
$spell.Range("J25").Value = "Damage 5"
$spell.Range("L25").Value = "Chorus"
$spell.Range("N25").Value = "If in control of Risen Undead deals 3 extra damage. "

# --- New column header -------------------------------------------------
$spell.Range("H1").Value = "Element Type"

# --- New "Element Type" column values (order matters for shared-string
# table layout: Shadow, Darkness, Water, Lightning, Earth, Ice are brand
# new strings; Time/Fire/Darkness(dup) reuse existing ones) ----------------
$spell.Range("H13").Value = "Shadow"
$spell.Range("H11").Value = "Darkness"
$spell.Range("H18").Value = "Water"
$spell.Range("H15").Value = "Lightning"
$spell.Range("H22").Value = "Earth"
$spell.Range("H25").Value = "Ice"
$spell.Range("H3").Value = "Time"
$spell.Range("H5").Value = "Time"
$spell.Range("H7").Value = "Fire"
$spell.Range("H9").Value = "Fire"
$spell.Range("H20").Value = "Darkness"

# --- Shift old Effect column (H) into new Effect column (J) ---------------
$spell.Range("J1").Value = "Effect"
$spell.Range("J3").Value = "15 damage"
$spell.Range("J5").Value = "Heal 20"
$spell.Range("J7").Value = "13 damage"
$spell.Range("J9").Value = "17 damage"
$spell.Range("J11").Value = "14 damage"
$spell.Range("J13").Value = "Defence 2"
$spell.Range("J15").Value = "Damage 10"
$spell.Range("J18").Value = "Damage 9"
$spell.Range("J20").Value = "Death 1"
$spell.Range("J22").Value = "Consume 1"

# --- Shift old Extra column (J) into new Extra column (L) -----------------
$spell.Range("L1").Value = "Extra"
$spell.Range("L3").Value = "None"
$spell.Range("L5").Value = "None"
$spell.Range("L7").Value = "Scorch"
$spell.Range("L9").Value = "None"
$spell.Range("L11").Value = "Stolen Life"
$spell.Range("L13").Value = "Shatter"
$spell.Range("L15").Value = "Energize"
$spell.Range("L16").Value = "Paralysis"
$spell.Range("L18").Value = "Magnifies"
$spell.Range("L20").Value = "None"
$spell.Range("L22").Value = "Unstable 6"

# --- Shift old Note column (L) into new Note column (N) -------------------
$spell.Range("N1").Value = "Note"
$spell.Range("N3").Value = "No extra"
$spell.Range("N5").Value = "Three turn cooldown"
$spell.Range("N7").Value = "deals 4 extra fire damge to weak opponets"
$spell.Range("N9").Value = "No extra"
$spell.Range("N11").Value = "Heal for 1/4 of total damage dealt"
$spell.Range("N13").Value = "After 3 turns Approaching Shadow will shatter dealing 32% of damage taken while it was active to enemy. Caps at 200."
$spell.Range("N15").Value = "Next physical attacks deal 3 Air damage. Next spell deals 5 Air damage, If spell is storm type deals extra 10 damage instead."
$spell.Range("N16").Value = "Has a .01% chance of paralyzing enemy for their next turn."
$spell.Range("N18").Value = "Deals 9 damage for 3 turns after which spell triggers dealing extra 30 damage. While spell charges all storm magic is amplified by 1 tier once triggered magnification ends."
$spell.Range("N20").Value = "Instantly kills enemy with 3 health or less. Upon death, if enemy is level 5 or less and 5 turns haven't passed since casting, collect soul. (soul counter +1)"
$spell.Range("N22").Value = "Raises Skeleton with 30 health, 20 physical damage, and no defence. Using spell consumes collected soul (Soul Counter -1)."
$spell.Range("N23").Value = " Each turn reduce unstable counter by 1, If counter equals 0 kill skeleton. If no soul counters are held spell fizzles."

# --- Clear the now-vacated old column positions ----------------------------
$spell.Range("J16").ClearContents()
$spell.Range("L23").ClearContents()

# --- Physical Armour: reword the magic-damage note (3% -> 20%) ------------
$phys.Range("P12").Value = "When hit by Magic Damage increase Phy Def by 2 and reduce speed by 20%."

# --- Selection / active tab bookkeeping ------------------------------------
$spell.Range("I27").Select()
$phys.Activate()
$phys.Range("P12").Select()
